$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column E (Order Status, Invoice Status)
$ws.Range("E:F").Insert()

$ws.Range("E:F").ColumnWidth = $ws.Range("D:D").ColumnWidth

$ws.Range("E4").Value = "Order Status"
$ws.Range("F4").Value = "Invoice Status"

$ws.Range("F5").Select()
